$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-enter the totals row formula across the range so Excel stores it as a
# single shared formula (t="shared") instead of one independent formula
# per cell.
$ws.Range("D7:W7").Formula = "=SUM(D4:D6)"

# Remove the blank spacer row (old row 24) and the extra blank row (old
# row 25) above the "Transformers / sub-satations" section, shifting that
# whole section up by two rows.
$ws.Rows("24:25").Delete() | Out-Null

# Leave the selection where Excel would land right after deleting the two
# selected rows.
$ws.Range("A24:XFD25").Select() | Out-Null
